# "rto month state  company script"
#
# The underlying data grid lists insurance rate-card rows keyed by
# company/product/vehicle-type/etc. Two data-entry corrections are made
# on the "bajaj" rows:
#   - Row 7 (Motor Insurance / Private Car Package Policy): the
#     vehical_subtype cell D7 had the wrong value copied in from a
#     different product row ("Non Petrol- Diesel, CNG,LPG"); it should
#     read "Bike 150 CC To 600 CC" instead, and gets the same highlighted
#     (boxed / shaded) header-style formatting used for the other
#     vehical_subtype / product_name header band cells in row 2.
#   - Row 8 (Life Insurance / Term Life Policy): the company name in A8
#     was abbreviated ("bajaj"); it is corrected to the full name
#     "Bajaj life line".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give D7 the same look (font/fill/border/alignment) as the other boxed
# header-style cells (C2/D2/E2) before overwriting its text.
$ws.Range("C2").Copy()
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("D7").Value = "Bike 150 CC To 600 CC"

# Correct the company name on row 8.
$ws.Range("A8").Value = "Bajaj life line"

# Column D (vehical_subtype) is now wide enough to show the corrected text.
$ws.Columns("D").ColumnWidth = 27

# Leave the selection on the cell that was just edited.
$ws.Range("D7").Select()
